# Update the "Share of Cost Effective Capacity Built in a Single Year" CSC
# sheet: recalibrate every resource's annual share value from 0.5 to 0.3.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

$ws.Range("B2:AE25").Value = 0.3
